{"js": "// Update the two-digit-by-two-digit multiplication \"answers\" worksheet:\n// each cell in the table holds a single equation like \"16\u00d737=592\" and\n// every one of the 25 populated cells gets replaced with a new equation,\n// per the commit's regenerated answer set. Run properties / formatting\n// on the existing run are preserved because we replace the text of the\n// matched (found) range in place rather than deleting + inserting a\n// brand-new run.\nconst replacements = [\n  [\"16\u00d737=592\", \"98\u00d797=9506\"],\n  [\"85\u00d759=5015\", \"61\u00d733=2013\"],\n  [\"52\u00d732=1664\", \"63\u00d715=945\"],\n  [\"11\u00d711=121\", \"67\u00d791=6097\"],\n  [\"52\u00d792=4784\", \"15\u00d733=495\"],\n  [\"24\u00d765=1560\", \"24\u00d775=1800\"],\n  [\"35\u00d726=910\", \"84\u00d746=3864\"],\n  [\"95\u00d736=3420\", \"11\u00d754=594\"],\n  [\"99\u00d794=9306\", \"18\u00d723=414\"],\n  [\"22\u00d773=1606\", \"55\u00d767=3685\"],\n  [\"30\u00d724=720\", \"91\u00d766=6006\"],\n  [\"86\u00d768=5848\", \"87\u00d736=3132\"],\n  [\"33\u00d764=2112\", \"81\u00d788=7128\"],\n  [\"76\u00d741=3116\", \"49\u00d727=1323\"],\n  [\"45\u00d774=3330\", \"95\u00d752=4940\"],\n  [\"88\u00d739=3432\", \"43\u00d795=4085\"],\n  [\"56\u00d746=2576\", \"21\u00d764=1344\"],\n  [\"84\u00d782=6888\", \"40\u00d795=3800\"],\n  [\"42\u00d794=3948\", \"73\u00d760=4380\"],\n  [\"62\u00d793=5766\", \"96\u00d723=2208\"],\n  [\"46\u00d712=552\", \"14\u00d734=476\"],\n  [\"65\u00d772=4680\", \"40\u00d784=3360\"],\n  [\"71\u00d750=3550\", \"98\u00d735=3430\"],\n  [\"39\u00d773=2847\", \"33\u00d780=2640\"],\n  [\"81\u00d748=3888\", \"59\u00d737=2183\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-by-two-digit multiplication \"answers\" worksheet:\n# each cell in the table holds a single equation like \"16x37=592\" and\n# every one of the 25 populated cells gets replaced with a new equation,\n# per the commit's regenerated answer set. Find/Replace on the document\n# Range preserves the existing run's formatting (font/size) since only\n# the text of the matched range is swapped.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old = \"16\u00d737=592\";   new = \"98\u00d797=9506\"},\n    @{old = \"85\u00d759=5015\";  new = \"61\u00d733=2013\"},\n    @{old = \"52\u00d732=1664\";  new = \"63\u00d715=945\"},\n    @{old = \"11\u00d711=121\";   new = \"67\u00d791=6097\"},\n    @{old = \"52\u00d792=4784\";  new = \"15\u00d733=495\"},\n    @{old = \"24\u00d765=1560\";  new = \"24\u00d775=1800\"},\n    @{old = \"35\u00d726=910\";   new = \"84\u00d746=3864\"},\n    @{old = \"95\u00d736=3420\";  new = \"11\u00d754=594\"},\n    @{old = \"99\u00d794=9306\";  new = \"18\u00d723=414\"},\n    @{old = \"22\u00d773=1606\";  new = \"55\u00d767=3685\"},\n    @{old = \"30\u00d724=720\";   new = \"91\u00d766=6006\"},\n    @{old = \"86\u00d768=5848\";  new = \"87\u00d736=3132\"},\n    @{old = \"33\u00d764=2112\";  new = \"81\u00d788=7128\"},\n    @{old = \"76\u00d741=3116\";  new = \"49\u00d727=1323\"},\n    @{old = \"45\u00d774=3330\";  new = \"95\u00d752=4940\"},\n    @{old = \"88\u00d739=3432\";  new = \"43\u00d795=4085\"},\n    @{old = \"56\u00d746=2576\";  new = \"21\u00d764=1344\"},\n    @{old = \"84\u00d782=6888\";  new = \"40\u00d795=3800\"},\n    @{old = \"42\u00d794=3948\";  new = \"73\u00d760=4380\"},\n    @{old = \"62\u00d793=5766\";  new = \"96\u00d723=2208\"},\n    @{old = \"46\u00d712=552\";   new = \"14\u00d734=476\"},\n    @{old = \"65\u00d772=4680\";  new = \"40\u00d784=3360\"},\n    @{old = \"71\u00d750=3550\";  new = \"98\u00d735=3430\"},\n    @{old = \"39\u00d773=2847\";  new = \"33\u00d780=2640\"},\n    @{old = \"81\u00d748=3888\";  new = \"59\u00d737=2183\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.Text = $p.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
